# "Homepage div adjusted and group member name added to the presentation"
#
# 1. Move the "TEAM MEMBERS" slide (last slide) up to position 2, right
#    after the title slide (Homepage div adjusted / re-ordering).
# 2. Append a new team member name ("OMOKEJIMI SAMUEL") as a new
#    paragraph in that slide's member-list placeholder.

$p = $ppt.ActivePresentation

# --- 1. Reorder slides -------------------------------------------------
$lastIndex = $p.Slides.Count
$p.Slides.Item($lastIndex).MoveTo(2)

# --- 2. Add the new team member's name ---------------------------------
# Re-fetch the slide by its new (post-move) index rather than reusing the
# old slide reference, since the handle tracks positional index.
$teamSlide = $p.Slides.Item(2)
$namesShape = $teamSlide.Shapes.Item(2)
$tr = $namesShape.TextFrame.TextRange
[void]$tr.InsertAfter("`rOMOKEJIMI SAMUEL")
